$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray " " value in F19 (was previously marked with a blank
# string placeholder) and mark Obadiah (row 19) as done, matching the
# "Edited until Obadiah 1:1" progress update.
$ws.Range("F19").Value = 1

# Update the current selection to reflect where editing left off.
$ws.Range("E19").Select()
